$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.06"
$ws.Range("E2").Value = "'5.82%"
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3").Value = "'32.49"
$ws.Range("E3").Value = "'11.40%"
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("D4").Value = "'5.302"
$ws.Range("E4").Value = "'1.79%"
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07504"
$ws.Range("E5").Value = "'7.45%"
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("D6").Value = "'7.835"
$ws.Range("E6").Value = "'5.84%"
$ws.Range("D6:E6").Style = "Normal"

$ws.Range("D7").Value = "'3.805"
$ws.Range("E7").Value = "'7.02%"
$ws.Range("D7:E7").Style = "Normal"

$ws.Range("D8").Value = "'1.497"
$ws.Range("E8").Value = "'6.44%"
$ws.Range("D8:E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9197"
$ws.Range("E9").Value = "'2.82%"
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'0.01776"
$ws.Range("E10").Value = "'2,644.73%"
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.1697"
$ws.Range("E11").Value = "'5.29%"
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07838"
$ws.Range("E12").Value = "'3.10%"
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("D13").Value = "'0.08042"
$ws.Range("E13").Value = "'5.12%"
$ws.Range("D13:E13").Style = "Normal"

$ws.Range("D14").Value = "'0.02995"
$ws.Range("E14").Value = "'2.48%"
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15").Value = "'0.09917"
$ws.Range("E15").Value = "'10.33%"
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'0.001493"
$ws.Range("E16").Value = "'-5.95%"
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'0.04621"
$ws.Range("E17").Value = "'2.33%"
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18").Value = "'0.006148"
$ws.Range("E18").Value = "'-4.96%"
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("D19").Value = "'3.469"
$ws.Range("E19").Value = "'0.26%"
$ws.Range("D19:E19").Style = "Normal"

$ws.Range("E20").Value = "'0.03%"
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'3.78%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.1334"
$ws.Range("E22").Value = "'0.54%"
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("D23").Value = "'4.499"
$ws.Range("E23").Value = "'12.34%"
$ws.Range("D23:E23").Style = "Normal"

$ws.Range("E25").Value = "'1.11%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.004450"
$ws.Range("E26").Value = "'1.99%"
$ws.Range("D26:E26").Style = "Normal"

$ws.Range("E27").Value = "'19.97%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.0001739"
$ws.Range("E28").Value = "'7.65%"
$ws.Range("D28:E28").Style = "Normal"

$ws.Range("D40").Value = "'0.04546"
$ws.Range("E40").Value = "'4.91%"
$ws.Range("D40:E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007181"
$ws.Range("E41").Value = "'3.64%"
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'8.46%"
$ws.Range("D42:E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002188"
$ws.Range("E43").Value = "'5.57%"
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01271"
$ws.Range("E44").Value = "'8.89%"
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006230"
$ws.Range("E45").Value = "'6.86%"
$ws.Range("D45:E45").Style = "Normal"

$ws.Range("D47").Value = "'0.01298"
$ws.Range("E47").Value = "'-0.49%"
$ws.Range("D47:E47").Style = "Normal"
